$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: conversion-of-the-day text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.89 = 40642.83 pesos`n✅ 40642.83 pesos = 9.86 = 956.88 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $text

# --- Update the "tasas" sheet: rate cells N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 101.1
$ws2.Range("O10").Value = 4108.99
$ws2.Range("N12").Value = 4120
$ws2.Range("O12").Value = 97
